$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07463599999999999
$ws.Range("N2").Value = 0.223908
$ws.Range("O2").Value = 0.06153261015633781
$ws.Range("P2").Value = 0.06153261015633781
$ws.Range("Q2").Value = 6.193868559115999
$ws.Range("R2").Value = 55.744817032044
$ws.Range("S2").Value = 0.02762509065299024
$ws.Range("T2").Value = 0.02762509065299024

# Row 3
$ws.Range("G3").Value = 82.98768099999999
$ws.Range("H3").Value = 248.963043
$ws.Range("I3").Value = 0.4489504115427952
$ws.Range("J3").Value = 0.4489504115427952
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.019266
$ws.Range("N3").Value = 0.057798
$ws.Range("O3").Value = 0.01588358523061263
$ws.Range("P3").Value = 0.01588358523061263
$ws.Range("Q3").Value = 1.598840662146
$ws.Range("R3").Value = 14.389565959314
$ws.Range("S3").Value = 0.007130942126058606
$ws.Range("T3").Value = 0.007130942126058606

# Row 4
$ws.Range("G4").Value = 82.98768099999999
$ws.Range("H4").Value = 248.963043
$ws.Range("I4").Value = 0.4489504115427952
$ws.Range("J4").Value = 0.4489504115427952
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.119048333333333
$ws.Range("N4").Value = 3.357145
$ws.Range("O4").Value = 0.9225838046130496
$ws.Range("P4").Value = 0.9225838046130496
$ws.Range("Q4").Value = 92.86722611024834
$ws.Range("R4").Value = 835.8050349922351
$ws.Range("S4").Value = 0.4141943787637464
$ws.Range("T4").Value = 0.4141943787637464

# Row 5
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("H5").Value = 189.421768
$ws.Range("I5").Value = 0.3415807409566563
$ws.Range("J5").Value = 0.3415807409566563
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.07463599999999999
$ws.Range("N5").Value = 0.223908
$ws.Range("O5").Value = 0.06153261015633781
$ws.Range("P5").Value = 0.06153261015633781
$ws.Range("Q5").Value = 4.712561025482666
$ws.Range("R5").Value = 42.413049229344
$ws.Range("S5").Value = 0.02101835457019894
$ws.Range("T5").Value = 0.02101835457019895

# Row 6
$ws.Range("G6").Value = 63.14058933333333
$ws.Range("H6").Value = 189.421768
$ws.Range("I6").Value = 0.3415807409566563
$ws.Range("J6").Value = 0.3415807409566563
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.019266
$ws.Range("N6").Value = 0.057798
$ws.Range("O6").Value = 0.01588358523061263
$ws.Range("P6").Value = 0.01588358523061263
$ws.Range("Q6").Value = 1.216466594096
$ws.Range("R6").Value = 10.948199346864
$ws.Range("S6").Value = 0.005425526812120865
$ws.Range("T6").Value = 0.005425526812120865

# Row 7
$ws.Range("G7").Value = 63.14058933333333
$ws.Range("H7").Value = 189.421768
$ws.Range("I7").Value = 0.3415807409566563
$ws.Range("J7").Value = 0.3415807409566563
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.119048333333333
$ws.Range("N7").Value = 3.357145
$ws.Range("O7").Value = 0.9225838046130496
$ws.Range("P7").Value = 0.9225838046130496
$ws.Range("Q7").Value = 70.65737125915112
$ws.Range("R7").Value = 635.9163413323599
$ws.Range("S7").Value = 0.3151368595743365
$ws.Range("T7").Value = 0.3151368595743365

# Row 8
$ws.Range("G8").Value = 38.719942
$ws.Range("H8").Value = 116.159826
$ws.Range("I8").Value = 0.2094688475005485
$ws.Range("J8").Value = 0.2094688475005485
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.07463599999999999
$ws.Range("N8").Value = 0.223908
$ws.Range("O8").Value = 0.06153261015633781
$ws.Range("P8").Value = 0.06153261015633781
$ws.Range("Q8").Value = 2.889901591112
$ws.Range("R8").Value = 26.009114320008
$ws.Range("S8").Value = 0.01288916493314862
$ws.Range("T8").Value = 0.01288916493314863

# Row 9
$ws.Range("G9").Value = 38.719942
$ws.Range("H9").Value = 116.159826
$ws.Range("I9").Value = 0.2094688475005485
$ws.Range("J9").Value = 0.2094688475005485
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.019266
$ws.Range("N9").Value = 0.057798
$ws.Range("O9").Value = 0.01588358523061263
$ws.Range("P9").Value = 0.01588358523061263
$ws.Range("Q9").Value = 0.7459784025720001
$ws.Range("R9").Value = 6.713805623148001
$ws.Range("S9").Value = 0.003327116292433162
$ws.Range("T9").Value = 0.003327116292433162

# Row 10
$ws.Range("G10").Value = 38.719942
$ws.Range("H10").Value = 116.159826
$ws.Range("I10").Value = 0.2094688475005485
$ws.Range("J10").Value = 0.2094688475005485
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.119048333333333
$ws.Range("N10").Value = 3.357145
$ws.Range("O10").Value = 0.9225838046130496
$ws.Range("P10").Value = 0.9225838046130496
$ws.Range("Q10").Value = 43.32948656186334
$ws.Range("R10").Value = 389.96537905677
$ws.Range("S10").Value = 0.1932525662749667
$ws.Range("T10").Value = 0.1932525662749667

